# Applies the cryptos-list price/volume/coin-row updates described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a literal text value even when the text looks numeric
# (e.g. "211.40"), preserving it as text (inlineStr) instead of letting Excel
# coerce it to a number and drop formatting such as trailing zeros.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "88.981.54"
$ws.Range("E2").Value = "  +0.70%  "

# Row 3
$ws.Range("D3").Value = "3.157.23"
$ws.Range("E3").Value = "  -3.96%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
Set-TextValue "D5" "211.40"
$ws.Range("E5").Value = "  -0.82%  "

# Row 6
Set-TextValue "D6" "612.09"
$ws.Range("E6").Value = "  -2.68%  "

# Row 7
Set-TextValue "D7" "0.384"
$ws.Range("E7").Value = "  +1.52%  "

# Row 8
Set-TextValue "D8" "0.685"
$ws.Range("E8").Value = "  -4.97%  "

# Row 9
$ws.Range("E9").Value = "  -0.04%  "

# Row 10
$ws.Range("D10").Value = "3.154.14"
$ws.Range("E10").Value = "  -3.92%  "

# Row 11
Set-TextValue "D11" "0.571"
$ws.Range("E11").Value = "  -1.33%  "

# Row 12
$ws.Range("E12").Value = "  -5.99%  "

# Row 13
$ws.Range("E13").Value = "  -5.06%  "

# Row 14
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "3.742.68"
$ws.Range("E14").Value = "  -3.78%  "

# Row 15
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "88.817.72"
$ws.Range("E15").Value = "  +0.60%  "

# Row 16
$ws.Range("B16").Value = "Toncoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D16" "5.21"
$ws.Range("E16").Value = "  -5.69%  "

# Row 17
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D17" "32.49"
$ws.Range("E17").Value = "  -5.65%  "

# Row 18
$ws.Range("D18").Value = "3.155.68"
$ws.Range("E18").Value = "  -4.71%  "

# Row 19
$ws.Range("E19").Value = "  +2.60%  "

# Row 20
Set-TextValue "D20" "13.31"
$ws.Range("E20").Value = "  -5.65%  "

# Row 21
Set-TextValue "D21" "433.72"
$ws.Range("E21").Value = "  -1.18%  "

# Row 22
$ws.Range("B22").Value = "PEPE"
$ws.Range("C22").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D22" "0.0000186"
$ws.Range("E22").Value = "  +36.52%  "

# Row 23
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D23" "8.54"
$ws.Range("E23").Value = "  -4.42%  "

# Row 24
Set-TextValue "D24" "5.05"
$ws.Range("E24").Value = "  -5.96%  "

# Row 25
Set-TextValue "D25" "5.08"
$ws.Range("E25").Value = "  -3.65%  "

# Row 26
Set-TextValue "D26" "11.65"
$ws.Range("E26").Value = "  -5.78%  "

# Row 27
$ws.Range("D27").Value = "3.327.23"
$ws.Range("E27").Value = "  -3.76%  "

# Row 28
Set-TextValue "D28" "74.79"
$ws.Range("E28").Value = "  -3.22%  "

# Row 29
$ws.Range("E29").Value = "  +0.00%  "

# Row 30
$ws.Range("E30").Value = "  -7.42%  "

# Row 31
$ws.Range("E31").Value = "  +0.28%  "

# Row 32
$ws.Range("E32").Value = "  +30.78%  "

# Row 33
Set-TextValue "D33" "8.37"
$ws.Range("E33").Value = "  -5.03%  "

# Row 34
Set-TextValue "D34" "527.83"
$ws.Range("E34").Value = "  -8.39%  "

# Row 35
Set-TextValue "D35" "6.97"
$ws.Range("E35").Value = "  -2.96%  "

# Row 36
$ws.Range("E36").Value = "  -6.46%  "

# Row 37
Set-TextValue "D37" "1.27"
$ws.Range("E37").Value = "  -8.75%  "

# Row 38
Set-TextValue "D38" "22.28"
$ws.Range("E38").Value = "  +2.05%  "

# Row 39
Set-TextValue "D39" "21.76"
$ws.Range("E39").Value = "  -4.82%  "

# Row 40
Set-TextValue "D40" "0.998"
$ws.Range("E40").Value = "  -0.04%  "

# Row 41
$ws.Range("E41").Value = "  -10.49%  "

# Row 42
$ws.Range("E42").Value = "  -0.02%  "

# Row 43
Set-TextValue "D43" "1.90"
$ws.Range("E43").Value = "  -6.43%  "

# Row 44
$ws.Range("E44").Value = "  -8.30%  "

# Row 45
Set-TextValue "D45" "149.48"
$ws.Range("E45").Value = "  -2.70%  "

# Row 46
Set-TextValue "D46" "43.68"
$ws.Range("E46").Value = "  -2.55%  "

# Row 47
Set-TextValue "D47" "170.96"
$ws.Range("E47").Value = "  -5.55%  "

# Row 48
$ws.Range("E48").Value = "  -10.69%  "

# Row 49
$ws.Range("E49").Value = "  -7.66%  "

# Row 50
Set-TextValue "D50" "4.02"
$ws.Range("E50").Value = "  -5.08%  "

# Row 51
Set-TextValue "D51" "0.604"
$ws.Range("E51").Value = "  -4.20%  "
